$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.176.66"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.403.78"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("D9").Value = "2.412.54"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.41"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D17").Value = "60.907.77"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "2.412.08"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  +5.01%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.17"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "594.48"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "0.0₃0947"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.372"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152.25"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.34"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.53"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.70"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.46"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  +1.37%  "
